# Apply cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.412.69"
$ws.Range("E2").Value = "  -1.69%  "

# Row 3
$ws.Range("D3").Value = "2.761.97"
$ws.Range("E3").Value = "  -2.75%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.91%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.19%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.551"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.34%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.22%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.49%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.137"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.78%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0836"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.90%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.51"
$ws.Range("D13").Style = "Normal"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.93%  "

# Row 15
$ws.Range("D15").Value = "3.199.86"
$ws.Range("E15").Value = "  -2.73%  "

# Row 16
$ws.Range("D16").Value = "2.760.73"
$ws.Range("E16").Value = "  -1.36%  "

# Row 17
$ws.Range("E17").Value = "  -1.79%  "

# Row 18
$ws.Range("D18").Value = "51.397.17"
$ws.Range("E18").Value = "  -1.53%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.75%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.36%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0960"
$ws.Range("E22").Value = "  -4.01%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.49%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.45%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.59%  "

# Row 26
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.54%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.163"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +13.61%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.70%  "

# Row 30
$ws.Range("E30").Value = "  -0.48%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.81%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "51.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.97%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.53%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0435"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.62%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0832"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.41%  "

# Row 36
$ws.Range("E36").Value = "  -8.19%  "

# Row 37
$ws.Range("E37").Value = "  +0.14%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.51%  "

# Row 39
$ws.Range("E39").Value = "  -5.21%  "

# Row 40
$ws.Range("E40").Value = "  -5.53%  "

# Row 41
$ws.Range("E41").Value = "  -3.38%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.83%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.59%  "

# Row 44
$ws.Range("E44").Value = "  -3.41%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.74%  "

# Row 46
$ws.Range("D46").Value = "2.077.83"
$ws.Range("E46").Value = "  -0.64%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.38%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.99%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.69%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.917"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.28%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.22%  "
